$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for every data row
# (rows 2-82). The sheet was refreshed, advancing that date by one day
# (46061 -> 46062) for every row.
for ($r = 2; $r -le 82; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value2 = $cell.Value2 + 1
}
